$d = $word.ActiveDocument

# Helper: build a minimal WordOpenXML single-package-part payload wrapping one <w:p>.
function New-ParaXml([string]$innerXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# 1. Title paragraph: consolidate "Ionic" / " Tutorial" runs (dropping the
#    spell-check proofErr wrapper) into a single run, and move the
#    _GoBack bookmark here (ahead of the run).
$titlePara = $d.Paragraphs.Item(1)
$titleXml = New-ParaXml('<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Ionic Tutorial</w:t></w:r>')
$titlePara.Range.InsertXML($titleXml) | Out-Null

$tbl = $d.Tables.Item(1)

# 2. Row 1 ("Criar aplicação"): consolidate the command cell text.
$cell1 = $tbl.Rows.Item(1).Cells.Item(2)
$cell1.Range.Paragraphs.Item(1).Range.InsertXML((New-ParaXml('<w:r><w:t>Ionic start myapp(tabs,sidemenu)</w:t></w:r>'))) | Out-Null

# 3. Row 2 ("Testar aplicação"): consolidate the command cell text.
$cell2 = $tbl.Rows.Item(2).Cells.Item(2)
$cell2.Range.Paragraphs.Item(1).Range.InsertXML((New-ParaXml('<w:r><w:t>Ionic serve</w:t></w:r>'))) | Out-Null

# 4. Row 3 ("Criar pagina"): consolidate the command cell text (this also
#    drops the _GoBack bookmark that used to live here - it moved to the
#    title paragraph above).
$cell3 = $tbl.Rows.Item(3).Cells.Item(2)
$cell3.Range.Paragraphs.Item(1).Range.InsertXML((New-ParaXml('<w:r><w:t>Ionic generate page &lt;nome&gt;</w:t></w:r>'))) | Out-Null

# 5. New row: "Toda vez que baixar o github a app" / "Npm install".
$newRow = $tbl.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Toda vez que baixar o github a app"
$newRow.Cells.Item(2).Range.Text = "Npm install"
